# Refresh the cryptocurrency price / 1h-volume figures to the latest scrape,
# as captured by the "Updated cryptos list ... with GitHub Actions" automation
# run. Two coin pairs also swapped rank position (rows 33/34: EthereumClassic
# <-> Fetch.AI, rows 45/46: ONDO <-> Filecoin), so their Coin/Link/Price/Volume
# cells move together as whole rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe (written as a doubled quote inside a single-quoted
# PowerShell string literal, e.g. '''578.84' == the text  '578.84 ) forces
# Excel to store the value as text instead of silently converting
# plain-decimal-looking Price strings (e.g. "578.84") into a floating point
# number -- matching the original inline-string cell type.

# Row 2
$ws.Range("D2").Value = '62.988.75'
$ws.Range("E2").Value = '  -0.49%  '

# Row 3
$ws.Range("D3").Value = '3.453.80'
$ws.Range("E3").Value = '  -0.92%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").Value = '''578.84'
$ws.Range("E5").Value = '  -1.15%  '

# Row 6
$ws.Range("D6").Value = '''148.89'
$ws.Range("E6").Value = '  +0.81%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("E8").Value = '  +0.09%  '

# Row 9
$ws.Range("D9").Value = '''7.88'
$ws.Range("E9").Value = '  +2.55%  '

# Row 10
$ws.Range("E10").Value = '  -2.20%  '

# Row 11
$ws.Range("D11").Value = '''0.408'
$ws.Range("E11").Value = '  +2.31%  '

# Row 12
$ws.Range("D12").Value = '4.044.95'
$ws.Range("E12").Value = '  -0.88%  '

# Row 13
$ws.Range("E13").Value = '  +2.30%  '

# Row 14
$ws.Range("E14").Value = '  -4.27%  '

# Row 15
$ws.Range("D15").Value = '3.452.24'
$ws.Range("E15").Value = '  -1.18%  '

# Row 16
$ws.Range("E16").Value = '  -1.41%  '

# Row 17
$ws.Range("D17").Value = '63.047.84'
$ws.Range("E17").Value = '  -0.40%  '

# Row 18
$ws.Range("E18").Value = '  +2.56%  '

# Row 19
$ws.Range("D19").Value = '''14.55'
$ws.Range("E19").Value = '  +1.15%  '

# Row 20
$ws.Range("D20").Value = '''9.19'
$ws.Range("E20").Value = '  -2.49%  '

# Row 21
$ws.Range("D21").Value = '''388.07'
$ws.Range("E21").Value = '  -1.27%  '

# Row 22
$ws.Range("D22").Value = '''0.561'
$ws.Range("E22").Value = '  -0.82%  '

# Row 23
$ws.Range("D23").Value = '''74.65'
$ws.Range("E23").Value = '  -0.79%  '

# Row 24
$ws.Range("E24").Value = '  -0.01%  '

# Row 25
$ws.Range("D25").Value = '3.594.83'
$ws.Range("E25").Value = '  -0.90%  '

# Row 26
$ws.Range("E26").Value = '  -3.81%  '

# Row 27
$ws.Range("D27").Value = '''0.183'
$ws.Range("E27").Value = '  -0.89%  '

# Row 28
$ws.Range("E28").Value = '  -2.61%  '

# Row 29
$ws.Range("E29").Value = '  +0.06%  '

# Row 30
$ws.Range("D30").Value = '''8.06'
$ws.Range("E30").Value = '  -2.42%  '

# Row 31
$ws.Range("E31").Value = '  -1.99%  '

# Row 32
$ws.Range("E32").Value = '  +0.03%  '

# Row 33
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = '''1.34'
$ws.Range("E33").Value = '  -5.74%  '

# Row 34
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").Value = '''23.34'
$ws.Range("E34").Value = '  -2.11%  '

# Row 35
$ws.Range("D35").Value = '''1.63'
$ws.Range("E35").Value = '  +3.11%  '

# Row 36
$ws.Range("D36").Value = '''5.35'
$ws.Range("E36").Value = '  +0.16%  '

# Row 37
$ws.Range("D37").Value = '''31.96'
$ws.Range("E37").Value = '  -2.05%  '

# Row 38
$ws.Range("E38").Value = '  -2.06%  '

# Row 39
$ws.Range("D39").Value = '''170.17'
$ws.Range("E39").Value = '  -0.97%  '

# Row 40
$ws.Range("D40").Value = '3.489.44'
$ws.Range("E40").Value = '  -0.95%  '

# Row 41
$ws.Range("E41").Value = '  +0.64%  '

# Row 42
$ws.Range("E42").Value = '  -1.06%  '

# Row 43
$ws.Range("D43").Value = '''42.88'
$ws.Range("E43").Value = '  +1.02%  '

# Row 44
$ws.Range("E44").Value = '  -1.88%  '

# Row 45
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").Value = '''4.36'
$ws.Range("E45").Value = '  -2.98%  '

# Row 46
$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").Value = '''1.19'
$ws.Range("E46").Value = '  -2.24%  '

# Row 47
$ws.Range("D47").Value = '2.583.03'
$ws.Range("E47").Value = '  -1.19%  '

# Row 48
$ws.Range("D48").Value = '''2.31'
$ws.Range("E48").Value = '  -0.30%  '

# Row 49
$ws.Range("E49").Value = '  +1.87%  '

# Row 50
$ws.Range("D50").Value = '''22.71'
$ws.Range("E50").Value = '  -4.77%  '

# Row 51
$ws.Range("E51").Value = '  +0.04%  '
